$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the next empty row after the last used row (row 86 -> new row 87)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($newRow, 2).Value = "temperature"

# Column C value looks numeric ("25") but must be stored as text, matching
# the rest of the sheet (which uses inline/shared strings throughout).
# Temporarily force a text number format so Excel doesn't coerce it to a
# number, then clear the formatting again so no stray cell style lingers.
$cCell = $ws.Cells.Item($newRow, 3)
$cCell.NumberFormat = "@"
$cCell.Value = "25"
$cCell.ClearFormats()

$ws.Cells.Item($newRow, 4).Value = "N/A"
$ws.Cells.Item($newRow, 5).Value = "N/A"
$ws.Cells.Item($newRow, 6).Value = "N/A"
